$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date in column C for every data row (2-342)
#    from 45184 to 45186.
for ($r = 2; $r -le 342; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# 2) Add a friendly display-name (the "Beteckning" value from column A) as the
#    second argument of the HYPERLINK() formulas found in columns S, T, U, V,
#    W, X and Y for the first ten data rows (2-11). Only row 2 has a formula
#    in column U.
$cols = @{
    "S" = @{ Path = "artfynd";        Ext = "xlsx" }
    "T" = @{ Path = "kartor";         Ext = "png"  }
    "U" = @{ Path = "knärot";         Ext = "png"  }
    "V" = @{ Path = "klagomål";       Ext = "docx" }
    "W" = @{ Path = "klagomålsmail";  Ext = "docx" }
    "X" = @{ Path = "tillsyn";        Ext = "docx" }
    "Y" = @{ Path = "tillsynsmail";   Ext = "docx" }
}

for ($r = 2; $r -le 11; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    foreach ($col in $cols.Keys) {
        if ($col -eq "U" -and $r -ne 2) {
            continue
        }
        $info = $cols[$col]
        $url = "https://klasma.github.io/Logging_LAXA/" + $info.Path + "/" + $beteckning + "." + $info.Ext
        $addr = $col + $r
        $ws.Range($addr).Formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
    }
}
